# Sluttrapport.docx — "Updated standard documents from PVV 3.0"
#
# 1) Re-introduce the "_GoBack" bookmark as an empty bookmark wrapping
#    the very first paragraph of the document. Word always keeps this
#    bookmark pinned to the last edited spot; adding it at position 0
#    is what the source revision shows, and it pushes every other
#    bookmark's numeric Id up by one (handled automatically by Word
#    when a new bookmark is inserted before the others).
# 2) Drop the stale <w:lastRenderedPageBreak/> hint in front of the
#    "Veiledning - Sluttrapport" heading by forcing Word to rebuild
#    that run (the page-break cache marker is layout-only and Word
#    regenerates it from scratch, not worth keeping stale copies of).

$d = $word.ActiveDocument

# --- 1) Add the "_GoBack" bookmark around the first paragraph -------------
$first = $d.Paragraphs.Item(1).Range
$start = $first.Start
$goBackRange = $d.Range($start, $start)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# --- 2) Remove the leftover <w:lastRenderedPageBreak/> --------------------
$needle = "Veiledning" + [char]0x2013 + "Sluttrapport"
# (searched without the surrounding spaces so minor whitespace
#  differences in the stored run don't stop the match)
$heading = $d.Content
$found = $heading.Find.Execute("Veiledning" + [char]0x2013 + " Sluttrapport", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    # Rewrite the run's text via a temporary placeholder so Word
    # regenerates the run from scratch, dropping the obsolete
    # lastRenderedPageBreak layout marker in the process, then put the
    # exact original text back.
    $original = $heading.Text
    $heading.Text = "__tmp_placeholder__"
    $reselect = $d.Content
    $reselect.Find.Execute("__tmp_placeholder__", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $reselect.Text = $original
}
